$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 text value update
$ws.Range("A2").Value = "test"

# Numeric cell updates in row 2 (values passed as strings then cast to double
# to avoid scientific-notation literal parsing issues)
$updates = @{
    "C2" = "4.882525671848932e-10"
    "D2" = "1.948945300608404e-08"
    "E2" = "4.266064354056439e-08"
    "F2" = "1.503641456335117e-23"
    "G2" = "4.14785212796965e-16"
    "H2" = "3.085076953267862e-09"
    "I2" = "0.3432081812120504"
    "J2" = "0.0002918850928256864"
    "K2" = "0.005021528492396222"
    "M2" = "0.01989653038264127"
    "N2" = "99.60035181888371"
    "O2" = "0.0165680418617828"
    "P2" = "0.003877809967144886"
    "Q2" = "0.0006680254130363223"
    "R2" = "0.003810555369340733"
    "S2" = "0.002972939742998978"
    "T2" = "0.003323261841090969"
    "U2" = "8.351425005809637e-06"
    "V2" = "9.753927993712402e-07"
    "W2" = "4.545095045639709e-12"
    "X2" = "1.966365578850191e-11"
    "Y2" = "1.038618059218076e-14"
    "Z2" = "1.120616220610563e-11"
    "AA2" = "1.706382763108598e-12"
    "AB2" = "4.007783232026285e-14"
    "AC2" = "4.151160139975632e-13"
    "AD2" = "8.479442672775411e-16"
    "AE2" = "5.740980206667846e-17"
    "AF2" = "2.259611602459314e-19"
    "AG2" = "1.102380258615769e-20"
    "AH2" = "3.259067833665462e-21"
    "AI2" = "5.123090910751867e-22"
    "AJ2" = "2.653374822683749e-22"
    "AL2" = "2.517892660083575e-08"
    "AM2" = "3.169419273690885e-09"
    "AN2" = "8.138013816879621e-10"
    "AR2" = "1.594581625582532e-18"
    "AS2" = "1.37547663983301e-16"
    "AT2" = "2.263516647216082e-16"
    "AU2" = "1.27872592821012e-31"
    "AV2" = "2.003708955537397e-24"
    "AW2" = "2.596106687162256e-17"
    "AX2" = "1.031928749934173e-09"
    "AY2" = "2.902031198891207e-12"
    "AZ2" = "5.328692796673984e-10"
    "BB2" = "2.631907674837067e-07"
    "BC2" = "0.01955592143757101"
    "BD2" = "3.348790211376283e-05"
    "BE2" = "0.01651599305331482"
    "BF2" = "0.01953907607943031"
    "BG2" = "2.745320513151094"
    "BH2" = "16.62495759095641"
    "BI2" = "36.12516140142506"
    "BJ2" = "14.63689524746325"
    "BK2" = "0.03138646094919402"
    "BL2" = "0.01232710543486643"
    "BM2" = "0.01761364812368117"
    "BN2" = "0.01438162264086357"
    "BO2" = "2.015481949372147"
    "BP2" = "1.056021978688168"
    "BQ2" = "5.349963934546627"
    "BR2" = "0.5937555741714557"
    "BS2" = "13.66459508903578"
    "BT2" = "0.001033688102286519"
    "BU2" = "2.181897666080586"
    "BV2" = "0.07190811079087266"
    "BW2" = "0.003888464267994268"
    "BX2" = "0.04772072675060171"
    "BY2" = "1.847011091668082"
    "CA2" = "2.124371281006574"
    "CB2" = "0.5444471050709344"
    "CC2" = "0.2342150070725434"
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = [double]$updates[$ref]
}

